$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 37
$ws.Range("A37").Value = 1
$ws.Range("B37").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C37").Value = "Arica y Parinacota"
$ws.Range("D37").Value = 44448
$ws.Range("D37").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E37").Value = 15
$ws.Range("F37").Value = 100112012
$ws.Range("G37").Value = "Espinaca"
$ws.Range("H37").Value = "Sin especificar"
$ws.Range("I37").Value = "Primera"
$ws.Range("J37").Value = 200
$ws.Range("K37").Value = 1400
$ws.Range("L37").Value = 1500
$ws.Range("M37").Value = 1450
$ws.Range("N37").Value = "$/atado 2,5 a 3 kilos"
$ws.Range("O37").Value = "Región de Arica y Parinacota"
$ws.Range("P37").Value = 483
$ws.Range("Q37").Value = 3
$ws.Range("R37").Value = "Hortaliza"

# Row 38
$ws.Range("A38").Value = 1
$ws.Range("B38").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C38").Value = "Arica y Parinacota"
$ws.Range("D38").Value = 44448
$ws.Range("D38").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E38").Value = 15
$ws.Range("F38").Value = 100112012
$ws.Range("G38").Value = "Espinaca"
$ws.Range("H38").Value = "Sin especificar"
$ws.Range("I38").Value = "Segunda"
$ws.Range("J38").Value = 200
$ws.Range("K38").Value = 1000
$ws.Range("L38").Value = 1200
$ws.Range("M38").Value = 1100
$ws.Range("N38").Value = "$/atado 2,5 a 3 kilos"
$ws.Range("O38").Value = "Región de Arica y Parinacota"
$ws.Range("P38").Value = 367
$ws.Range("Q38").Value = 3
$ws.Range("R38").Value = "Hortaliza"
